# Update "want to go" counts (column F) on the "展览" and "全部类型" sheets
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 14796
    $ws.Range("F3").Value = 18258
}

# Row with F=1200 -> 1201 is row 26 on "展览" and row 27 on "全部类型"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F26").Value = 1201
$ws1.Range("F34").Value = 5242

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F27").Value = 1201
$ws4.Range("F36").Value = 5242
